$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Locations_mapping: insert a new row 2 with a "*DEFAULT*" -> "*DEFAULT*"
#    mapping entry above the existing data row.
# ---------------------------------------------------------------------
$loc = $wb.Worksheets.Item("Locations_mapping")
$loc.Range("A2:D2").Insert()
$loc.Range("A2").Value = "*DEFAULT*"
$loc.Range("B2").Value = "*DEFAULT*"
$loc.Range("A2:B2").Select()

# ---------------------------------------------------------------------
# 2. Add a new "Item_policies" worksheet, placed right before the
#    data_validation sheet, holding the source/destination item policy
#    code mapping.
# ---------------------------------------------------------------------
$dv = $wb.Worksheets.Item("data_validation")
$ip = $wb.Worksheets.Add($dv)
$ip.Name = "Item_policies"

# Column widths (characters).
$ip.Columns.Item(1).ColumnWidth = 22.1666666667
$ip.Columns.Item(2).ColumnWidth = 24.6080729167

# Pre-register the bold + text-number-format style ahead of the plain
# text-number-format one (matches the style order produced by the
# original authoring session) using a scratch cell that is discarded
# straight away.
$ip.Range("Z1").NumberFormat = "@"
$ip.Range("Z1").Font.Bold = $true
$ip.Range("Z1").Value = "tmp"
$ip.Columns.Item(26).Delete()

$ip.Range("A1").Value = "Source item policy code"
$ip.Range("B1").Value = "Destination item policy code"
$ip.Range("A2").Value = "*DEFAULT*"

# Keep "01" as text, not a number.
$ip.Range("B2").NumberFormat = "@"
$ip.Range("B2").Value = "01"

# Apply a text number format to both columns and bold the whole table.
$ip.Columns.Item(1).NumberFormat = "@"
$ip.Columns.Item(2).NumberFormat = "@"
$ip.Range("A1:B2").Font.Bold = $true

$ip.PageSetup.PaperSize = 9
$ip.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# 3. Hide the data_validation sheet (still used by data-validation
#    lists elsewhere, but no longer meant to be browsed directly).
#    NB: re-fetch the sheet by name - the handle captured before
#    Worksheets.Add() above tracks the newly inserted sheet instead.
# ---------------------------------------------------------------------
$dv2 = $wb.Worksheets.Item("data_validation")
$dv2.Visible = $false

# ---------------------------------------------------------------------
# 4. Leave the new sheet as the active one / active selection, matching
#    where the author was last working.
# ---------------------------------------------------------------------
$ip.Range("F10").Select()
